$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells (one at a time, since
# this runtime only applies NumberFormat to the first area of a union range)
# so Excel does not auto-convert them to numbers: they must stay literal text.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.681.43"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.844.04"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "315.46"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("D8").Value = "0.3709"
$ws.Range("E8").Value = "  +2.11%  "
$ws.Range("D9").Value = "0.07337"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").Value = "0.8779"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").Value = "21.05"
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").Value = "1.882.58"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "5.474"
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Value = "0.06960"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "81.13"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "0.000009052"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").Value = "0.9998"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").Value = "15.58"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").Value = "27.819.52"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").Value = "5.099"
$ws.Range("E22").Value = "  +2.81%  "
$ws.Range("D23").Value = "10.98"
$ws.Range("E23").Value = "  +5.67%  "
$ws.Range("D24").Value = "2.134.68"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").Value = "1.989"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "154.15"
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("D27").Value = "18.94"
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("D28").Value = "5.325"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").Value = "115.71"
$ws.Range("E29").Value = "  -4.76%  "
$ws.Range("D30").Value = "1.876"
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("D31").Value = "0.08925"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("D32").Value = "0.7900"
$ws.Range("E32").Value = "  +3.53%  "
$ws.Range("D33").Value = "4.619"
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("E34").Value = "  +6.39%  "
$ws.Range("D35").Value = "2.972"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "0.05444"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").Value = "1.105"
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("D39").Value = "0.01959"
$ws.Range("E39").Value = "  +1.38%  "
$ws.Range("D40").Value = "2.835"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").Value = "0.5175"
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("D42").Value = "0.1694"
$ws.Range("E42").Value = "  +2.48%  "
$ws.Range("D43").Value = "6.796"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").Value = "8.645"
$ws.Range("E44").Value = "  +3.28%  "
$ws.Range("D45").Value = "10.67"
$ws.Range("E45").Value = "  +3.59%  "
$ws.Range("D46").Value = "0.4794"
$ws.Range("E46").Value = "  +2.59%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.06557"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "106.52"
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("D49").Value = "0.9997"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "1.663"
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("D51").Value = "1.845"
$ws.Range("E51").Value = "  +5.67%  "
